$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: overwrite the existing employee record (Leandro Moraga Semir / 19 / 125000)
# with a new employee (Parbati Sano Poudel). The Age/Salary columns are cleared
# entirely for this record (no numeric data supplied for the new employee).
$ws.Range("A2").Value = "Parbati"
$ws.Range("B2").Value = "Sano"
$ws.Range("C2").Value = "Poudel"
$ws.Range("D2:E2").ClearContents()

# Append two more new employee records (same data) as rows 6 and 7.
$ws.Range("A6").Value = "Parbati"
$ws.Range("B6").Value = "Sano"
$ws.Range("C6").Value = "Poudel"

$ws.Range("A7").Value = "Parbati"
$ws.Range("B7").Value = "Sano"
$ws.Range("C7").Value = "Poudel"
